# Update cryptos list cell values per the latest data pull.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @(
    @('D2', '43.608.91'),
    @('E2', '  -0.60%  '),
    @('D3', '2.222.08'),
    @('E3', '  +0.60%  '),
    @('E4', '  -0.14%  '),
    @('D5', '''270.20'),
    @('E5', '  +4.01%  '),
    @('D6', '''92.78'),
    @('E6', '  +13.77%  '),
    @('D7', '''0.624'),
    @('E7', '  -0.59%  '),
    @('E8', '  -0.02%  '),
    @('D9', '''0.618'),
    @('E9', '  +3.05%  '),
    @('D10', '''45.76'),
    @('E10', '  +5.63%  '),
    @('D11', '''0.0935'),
    @('E11', '  +1.42%  '),
    @('D12', '''8.27'),
    @('E12', '  +18.04%  '),
    @('E13', '  +0.83%  '),
    @('D14', '2.555.84'),
    @('E14', '  +0.59%  '),
    @('E15', '  +4.09%  '),
    @('B16', 'Polygon'),
    @('C16', 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'),
    @('D16', '''0.801'),
    @('E16', '  +3.14%  '),
    @('B17', 'WrappedEther'),
    @('C17', 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'),
    @('D17', '2.218.94'),
    @('E17', '  -0.13%  '),
    @('D18', '43.524.97'),
    @('E18', '  -0.64%  '),
    @('E19', '  +0.94%  '),
    @('D20', '''6.03'),
    @('E20', '  +0.92%  '),
    @('D21', '''70.43'),
    @('E21', '  -0.88%  '),
    @('E22', '  -1.43%  '),
    @('D23', '''232.56'),
    @('E23', '  +0.35%  '),
    @('D24', '''9.12'),
    @('E24', '  -0.96%  '),
    @('E25', '  -0.01%  '),
    @('D26', '''11.36'),
    @('E26', '  +5.98%  '),
    @('D27', '''2.51'),
    @('E27', '  +11.51%  '),
    @('D28', '''41.81'),
    @('E28', '  +1.79%  '),
    @('E29', '  +5.19%  '),
    @('E30', '  +1.87%  '),
    @('D31', '''172.65'),
    @('E31', '  -0.06%  '),
    @('E32', '  +5.58%  '),
    @('D33', '''20.88'),
    @('E33', '  +1.83%  '),
    @('D34', '''5.48'),
    @('E34', '  +3.26%  '),
    @('E35', '  +0.82%  '),
    @('E36', '  -2.38%  '),
    @('E37', '  -2.09%  '),
    @('D38', '''4.30'),
    @('E38', '  -4.66%  '),
    @('D39', '''3.60'),
    @('E39', '  +23.77%  '),
    @('D40', '''12.58'),
    @('E40', '  -4.79%  '),
    @('D41', '''0.220'),
    @('E41', '  +9.94%  '),
    @('D42', '''2.16'),
    @('E42', '  +2.51%  '),
    @('D43', '''63.35'),
    @('E43', '  +0.92%  '),
    @('D44', '''5.32'),
    @('E44', '  -3.44%  '),
    @('E45', '  +0.11%  '),
    @('E46', '  +0.68%  '),
    @('D47', '''100.24'),
    @('E47', '  -1.84%  '),
    @('E48', '  +3.25%  '),
    @('E49', '  +1.76%  '),
    @('D50', '''0.438'),
    @('E50', '  -0.82%  '),
    @('E51', '  -5.04%  ')
)

foreach ($u in $updates) {
    $ws.Range($u[0]).Value = $u[1]
}
